$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.1762295081967213
$ws.Cells.Item(2, 3).Value = 0.5860655737704918
$ws.Cells.Item(2, 10).Value = 0.02049180327868852
$ws.Cells.Item(2, 16).Value = 0.1434426229508197
$ws.Cells.Item(2, 19).Value = 0.07377049180327869
$ws.Cells.Item(3, 3).Value = 0.05333333333333334
$ws.Cells.Item(3, 10).Value = 0.02
$ws.Cells.Item(3, 16).Value = 0.76
$ws.Cells.Item(3, 19).Value = 0.1666666666666667
$ws.Cells.Item(4, 10).Value = 0.06521739130434782
$ws.Cells.Item(4, 16).Value = 0.5869565217391305
$ws.Cells.Item(4, 19).Value = 0.3478260869565217
$ws.Cells.Item(6, 2).Value = 0.05627705627705628
$ws.Cells.Item(6, 4).Value = 0.01298701298701299
$ws.Cells.Item(6, 6).Value = 0.1082251082251082
$ws.Cells.Item(6, 10).Value = 0.2121212121212121
$ws.Cells.Item(6, 15).Value = 0.02164502164502164
$ws.Cells.Item(6, 17).Value = 0.1818181818181818
$ws.Cells.Item(6, 18).Value = 0.06493506493506493
$ws.Cells.Item(6, 19).Value = 0.341991341991342
$ws.Cells.Item(7, 2).Value = 0.08823529411764706
$ws.Cells.Item(7, 4).Value = 0.0196078431372549
$ws.Cells.Item(7, 5).Value = 0.004901960784313725
$ws.Cells.Item(7, 6).Value = 0.06862745098039216
$ws.Cells.Item(7, 10).Value = 0.107843137254902
$ws.Cells.Item(7, 15).Value = 0.05392156862745098
$ws.Cells.Item(7, 17).Value = 0.1911764705882353
$ws.Cells.Item(7, 18).Value = 0.06372549019607843
$ws.Cells.Item(7, 19).Value = 0.4019607843137255
$ws.Cells.Item(8, 2).Value = 0.07709750566893424
$ws.Cells.Item(8, 4).Value = 0.01587301587301587
$ws.Cells.Item(8, 6).Value = 0.04988662131519275
$ws.Cells.Item(8, 10).Value = 0.1133786848072562
$ws.Cells.Item(8, 15).Value = 0.0272108843537415
$ws.Cells.Item(8, 17).Value = 0.2018140589569161
$ws.Cells.Item(8, 18).Value = 0.09977324263038549
$ws.Cells.Item(8, 19).Value = 0.4149659863945578
$ws.Cells.Item(9, 2).Value = 0.08636363636363636
$ws.Cells.Item(9, 4).Value = 0.01818181818181818
$ws.Cells.Item(9, 6).Value = 0.08636363636363636
$ws.Cells.Item(9, 10).Value = 0.08181818181818182
$ws.Cells.Item(9, 15).Value = 0.04090909090909091
$ws.Cells.Item(9, 17).Value = 0.1409090909090909
$ws.Cells.Item(9, 18).Value = 0.09545454545454546
$ws.Cells.Item(9, 19).Value = 0.45
$ws.Cells.Item(10, 2).Value = 0.101063829787234
$ws.Cells.Item(10, 4).Value = 0.02482269503546099
$ws.Cells.Item(10, 5).Value = 0.0008865248226950354
$ws.Cells.Item(10, 6).Value = 0.07535460992907801
$ws.Cells.Item(10, 10).Value = 0.1312056737588652
$ws.Cells.Item(10, 15).Value = 0.02216312056737589
$ws.Cells.Item(10, 17).Value = 0.2349290780141844
$ws.Cells.Item(10, 18).Value = 0.07003546099290781
$ws.Cells.Item(10, 19).Value = 0.3395390070921986
$ws.Cells.Item(11, 7).Value = 0.1383399209486166
$ws.Cells.Item(11, 10).Value = 0.04347826086956522
$ws.Cells.Item(11, 11).Value = 0.1541501976284585
$ws.Cells.Item(11, 12).Value = 0.6521739130434783
$ws.Cells.Item(11, 19).Value = 0.01185770750988142
$ws.Cells.Item(12, 7).Value = 0.807909604519774
$ws.Cells.Item(12, 10).Value = 0.1299435028248588
$ws.Cells.Item(12, 11).Value = 0.005649717514124294
$ws.Cells.Item(12, 12).Value = 0.05084745762711865
$ws.Cells.Item(12, 19).Value = 0.005649717514124294
$ws.Cells.Item(13, 7).Value = 0.725
$ws.Cells.Item(13, 10).Value = 0.25
$ws.Cells.Item(13, 19).Value = 0.025
$ws.Cells.Item(15, 6).Value = 0.02164502164502164
$ws.Cells.Item(15, 8).Value = 0.1688311688311688
$ws.Cells.Item(15, 9).Value = 0.06493506493506493
$ws.Cells.Item(15, 10).Value = 0.3463203463203463
$ws.Cells.Item(15, 11).Value = 0.04329004329004329
$ws.Cells.Item(15, 13).Value = 0.004329004329004329
$ws.Cells.Item(15, 15).Value = 0.08225108225108226
$ws.Cells.Item(15, 19).Value = 0.2683982683982684
$ws.Cells.Item(16, 6).Value = 0.01162790697674419
$ws.Cells.Item(16, 8).Value = 0.2151162790697674
$ws.Cells.Item(16, 9).Value = 0.1162790697674419
$ws.Cells.Item(16, 10).Value = 0.3604651162790697
$ws.Cells.Item(16, 11).Value = 0.08139534883720931
$ws.Cells.Item(16, 13).Value = 0.02325581395348837
$ws.Cells.Item(16, 15).Value = 0.06976744186046512
$ws.Cells.Item(16, 19).Value = 0.1220930232558139
$ws.Cells.Item(17, 6).Value = 0.01079913606911447
$ws.Cells.Item(17, 8).Value = 0.1555075593952484
$ws.Cells.Item(17, 9).Value = 0.1101511879049676
$ws.Cells.Item(17, 10).Value = 0.4341252699784017
$ws.Cells.Item(17, 11).Value = 0.1187904967602592
$ws.Cells.Item(17, 13).Value = 0.01943844492440605
$ws.Cells.Item(17, 14).Value = 0.002159827213822894
$ws.Cells.Item(17, 15).Value = 0.06047516198704104
$ws.Cells.Item(17, 19).Value = 0.08855291576673865
$ws.Cells.Item(18, 6).Value = 0.04117647058823529
$ws.Cells.Item(18, 8).Value = 0.2352941176470588
$ws.Cells.Item(18, 9).Value = 0.08235294117647059
$ws.Cells.Item(18, 10).Value = 0.3647058823529412
$ws.Cells.Item(18, 11).Value = 0.07647058823529412
$ws.Cells.Item(18, 13).Value = 0.02352941176470588
$ws.Cells.Item(18, 14).Value = 0.005882352941176471
$ws.Cells.Item(18, 15).Value = 0.07058823529411765
$ws.Cells.Item(18, 19).Value = 0.1
$ws.Cells.Item(19, 6).Value = 0.01610017889087657
$ws.Cells.Item(19, 8).Value = 0.2209302325581395
$ws.Cells.Item(19, 9).Value = 0.1046511627906977
$ws.Cells.Item(19, 10).Value = 0.3694096601073345
$ws.Cells.Item(19, 11).Value = 0.1073345259391771
$ws.Cells.Item(19, 13).Value = 0.01967799642218247
$ws.Cells.Item(19, 15).Value = 0.07155635062611806
$ws.Cells.Item(19, 19).Value = 0.09033989266547406
